$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (rows 5-7), mirroring existing rows 2-4 but with
# updated values and a different source filename.
$newRows = @(
    @("§ 275.0-2_P1|llm_response", 9, 0, 7, 2, 20, 65, "documents-2024-11-01-1.json"),
    @("§ 275.0-5_P1|llm_response", 5, 0, 3, 2, 15, 25, "documents-2024-11-01-1.json"),
    @("§ 275.0-7_P1|llm_response", 9, 0, 5, 4, 19, 36, "documents-2024-11-01-1.json")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $ws.Cells.Item($r, 8).Value = $rowData[7]
}
